# Auto-generated PowerShell Excel COM-interop script
# Updates 'Price' (column D) and 'Volume(1h)' (column E) values in the
# cryptos worksheet to reflect freshly fetched market data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.901.77'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').Value = '3.511.24'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'608.32"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.07%  '
$ws.Range('D6').Value = "'198.62"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.75%  '
$ws.Range('E7').Value = '  +1.40%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').Value = "'0.211"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.71%  '
$ws.Range('E10').Value = '  +2.24%  '
$ws.Range('D11').Value = "'54.45"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.05%  '
$ws.Range('E12').Value = '  -0.81%  '
$ws.Range('D13').Value = "'9.63"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.19%  '
$ws.Range('D14').Value = '4.066.22'
$ws.Range('E14').Value = '  -1.37%  '
$ws.Range('D15').Value = "'597.88"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.91%  '
$ws.Range('D16').Value = '69.981.47'
$ws.Range('E16').Value = '  -0.40%  '
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('D19').Value = '3.503.97'
$ws.Range('E19').Value = '  -1.95%  '
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').Value = "'0.998"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.56%  '
$ws.Range('D22').Value = "'17.77"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.49%  '
$ws.Range('D23').Value = "'103.97"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +10.40%  '
$ws.Range('E24').Value = '  -0.94%  '
$ws.Range('E25').Value = '  +4.67%  '
$ws.Range('E26').Value = '  +6.66%  '
$ws.Range('D27').Value = "'11.02"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').Value = "'9.84"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.52%  '
$ws.Range('D29').Value = "'33.90"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.17%  '
$ws.Range('D30').Value = "'4.57"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +22.68%  '
$ws.Range('E31').Value = '  +3.13%  '
$ws.Range('D32').Value = "'12.80"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.76%  '
$ws.Range('E33').Value = '  +1.75%  '
$ws.Range('D34').Value = "'63.79"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('D35').Value = '3.712.88'
$ws.Range('E35').Value = '  +2.84%  '
$ws.Range('D36').Value = "'526.56"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.81%  '
$ws.Range('E37').Value = '  -0.03%  '
$ws.Range('D38').Value = '0.0₃0802'
$ws.Range('E38').Value = '  +3.07%  '
$ws.Range('D39').Value = "'3.02"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.86%  '
$ws.Range('E40').Value = '  -2.82%  '
$ws.Range('D41').Value = "'36.96"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.99%  '
$ws.Range('E42').Value = '  +1.59%  '
$ws.Range('E43').Value = '  +0.32%  '
$ws.Range('E44').Value = '  +1.61%  '
$ws.Range('D45').Value = "'2.88"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.47%  '
$ws.Range('D46').Value = "'0.140"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.97%  '
$ws.Range('E47').Value = '  -4.38%  '
$ws.Range('D48').Value = "'8.79"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E49').Value = '  +0.31%  '
$ws.Range('D50').Value = "'132.20"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.27%  '
$ws.Range('E51').Value = '  -1.86%  '
